$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (30) down to new row 31
$ws.Range("A30:D30").Copy()
$ws.Range("A31:D31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in new row of data (row 31)
$ws.Range("A31").Value = 43696
$ws.Range("B31").Value = 0.51944444444444449
$ws.Range("C31").Value = 55
$ws.Range("D31").Value = 19928

# Update the view: scroll down one row (topLeftCell A10 -> A11) and move
# selection to C32 (mirrors the user typing into the next row and Excel
# auto-scrolling by one row)
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C32").Select()
